$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSheet")

# Row 6: same values as row 4/5 (Abhi0, Abhi1, Abhi2)
$ws.Range("A6").Value = "Abhi0"
$ws.Range("B6").Value = "Abhi1"
$ws.Range("C6").Value = "Abhi2"

# Rows 7-11: new values (Abhi_0, Abhi_1, Abhi_2) repeated 5 times
for ($r = 7; $r -le 11; $r++) {
    $ws.Range("A$r").Value = "Abhi_0"
    $ws.Range("B$r").Value = "Abhi_1"
    $ws.Range("C$r").Value = "Abhi_2"
}
